$wb = $excel.ActiveWorkbook

# --- MasterDataCreation (sheet2): add a new "Join Name" column (AM) and
#     update the PV manual join timestamp value in AL2 ---
$wsMaster = $wb.Worksheets.Item("MasterDataCreation")
$wsMaster.Range("AL2").Value = "A_PVMJ_FriJan181634082019"
$wsMaster.Range("AM1").Value = "Join Name"
$wsMaster.Range("AM2").Value = "A_PVMJ_FriJan181641252019"

# match the header formatting (dark-blue fill / white font) used by the
# rest of row 1, e.g. the adjacent AL1 cell
$wsMaster.Range("AL1").Copy()
$wsMaster.Range("AM1").PasteSpecial(-4122)

# --- CongaTemplateCreation (sheet3): refresh the generated Name value ---
$wsConga = $wb.Worksheets.Item("CongaTemplateCreation")
$wsConga.Range("F2").Value = "AFriJan181534322019"
